$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.410.29"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.697.25"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.45"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5477"
$ws.Range("E6").Value = "  +3.88%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2739"
$ws.Range("E8").Value = "  +1.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06448"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.97"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07674"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").Value = "1.696.08"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.558"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5854"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008406"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.67"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "26.454.82"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.943"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.40"
$ws.Range("E21").Value = "  +0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.271"
$ws.Range("E22").Value = "  +0.89%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.96"
$ws.Range("E24").Value = "  +2.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1312"
$ws.Range("E25").Value = "  +5.27%  "
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.81"
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06234"
$ws.Range("E28").Value = "  -6.14%  "
$ws.Range("E29").Value = "  +2.27%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.611"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.686"
$ws.Range("E33").Value = "  +1.35%  "
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6160"
$ws.Range("E35").Value = "  -0.99%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.411"
$ws.Range("E36").Value = "  +0.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.767"
$ws.Range("E37").Value = "  +1.48%  "
$ws.Range("E38").Value = "  +1.98%  "
$ws.Range("D39").Value = "1.117.90"
$ws.Range("E39").Value = "  +0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.113"
$ws.Range("E40").Value = "  -5.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8828"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.19"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("D44").Value = "1.847.84"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000109"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.66"
$ws.Range("E46").Value = "  +1.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.217"
$ws.Range("E47").Value = "  +0.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05283"
$ws.Range("E49").Value = "  +0.22%  "
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4301"
$ws.Range("E51").Value = "  +0.01%  "
